$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-45: each entry is (row, I-value, J-value)
$rows = @(
    @(2,9,9),
    @(3,7,7),
    @(4,6,6),
    @(5,6,6),
    @(6,5,5),
    @(7,6,7),
    @(8,6,7),
    @(9,9,9),
    @(10,6,6),
    @(11,8,8),
    @(12,7,8),
    @(13,6,7),
    @(14,7,8),
    @(15,6,6),
    @(16,6,6),
    @(17,5,6),
    @(18,6,7),
    @(19,6,7),
    @(20,7,8),
    @(21,7,7),
    @(22,8,8),
    @(23,7,8),
    @(24,6,7),
    @(25,8,8),
    @(26,7,7),
    @(27,6,6),
    @(28,6,6),
    @(29,8,8),
    @(30,9,9),
    @(31,7,8),
    @(32,6,7),
    @(33,7,7),
    @(34,7,7),
    @(35,9,9),
    @(36,7,8),
    @(37,8,8),
    @(38,5,5),
    @(39,8,9),
    @(40,7,7),
    @(41,8,8),
    @(42,6,6),
    @(43,4,4),
    @(44,4,4),
    @(45,3,3)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value  = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

Write-Output "applied I0/IF columns"
